# Fruta / hortaliza, semanal
# Insert two new weekly records at the top of the data block (rows 39-40),
# pushing all the existing records (previously rows 39-135) down by two rows
# (to rows 41-137).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the first data row of the block (row 39).
# Excel will shift every row from 39 downward by two rows, automatically
# turning the former row 135 into row 137, and extending the used range.
$ws.Rows.Item(39).Resize(2).Insert()

# --- New record #1 -> row 39 ---
$ws.Cells.Item(39, 1).Value  = 5
$ws.Cells.Item(39, 2).Value  = 'Macroferia Regional de Talca'
$ws.Cells.Item(39, 3).Value  = 'Maule'
$ws.Cells.Item(39, 4).Value  = 44544
$ws.Cells.Item(39, 5).Value  = 7
$ws.Cells.Item(39, 6).Value  = 100112024
$ws.Cells.Item(39, 7).Value  = 'Choclo'
$ws.Cells.Item(39, 8).Value  = 'Choclero'
$ws.Cells.Item(39, 9).Value  = 'Primera'
$ws.Cells.Item(39, 10).Value = 20000
$ws.Cells.Item(39, 11).Value = 400
$ws.Cells.Item(39, 12).Value = 400
$ws.Cells.Item(39, 13).Value = 400
$ws.Cells.Item(39, 14).Value = '$/unidad'
$ws.Cells.Item(39, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(39, 16).Value = 400
$ws.Cells.Item(39, 17).Value = 1
$ws.Cells.Item(39, 18).Value = 'Hortaliza'

# --- New record #2 -> row 40 ---
$ws.Cells.Item(40, 1).Value  = 5
$ws.Cells.Item(40, 2).Value  = 'Macroferia Regional de Talca'
$ws.Cells.Item(40, 3).Value  = 'Maule'
$ws.Cells.Item(40, 4).Value  = 44544
$ws.Cells.Item(40, 5).Value  = 7
$ws.Cells.Item(40, 6).Value  = 100112024
$ws.Cells.Item(40, 7).Value  = 'Choclo'
$ws.Cells.Item(40, 8).Value  = 'Choclero'
$ws.Cells.Item(40, 9).Value  = 'Segunda'
$ws.Cells.Item(40, 10).Value = 20000
$ws.Cells.Item(40, 11).Value = 300
$ws.Cells.Item(40, 12).Value = 300
$ws.Cells.Item(40, 13).Value = 300
$ws.Cells.Item(40, 14).Value = '$/unidad'
$ws.Cells.Item(40, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(40, 16).Value = 300
$ws.Cells.Item(40, 17).Value = 1
$ws.Cells.Item(40, 18).Value = 'Hortaliza'
